$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product rows to append (rows 81-85).
$newRows = @(
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00"),
    @("9092831390928320", "Ar Condicionado Split 9000 BTUs Frio Air Volution Springer Midea", "1.743.15"),
    @("9092831390928320", "Ar Condicionado Split 9000 BTUs Frio Air Volution Springer Midea", "1.743.15")
)

$startRow = 81
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    # Column A holds purely numeric-looking codes (e.g. "91989296"). A plain
    # .Value assignment would be auto-coerced into a real number by Excel, so
    # instead compute it through a TEXT() formula and convert that formula to
    # a static value via copy / paste-special (values only). This keeps the
    # cell a genuine text string (matching the source inlineStr cells) without
    # leaving behind any extra/unused cell style definitions.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = '=TEXT(' + $rowData[0] + ',"0")'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    # Column B is free-form text that never looks like a number, so a direct
    # value assignment is stored as text already.
    $ws.Cells.Item($r, 2).Value = $rowData[1]

    # Column C values contain multiple "." separators (e.g. "8.999.00"),
    # which are not valid numbers, so Excel keeps them as text automatically.
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

$excel.CutCopyMode = 0
